$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '71.777.60'
$ws.Range("E2").Value = '  -1.39%  '

# Row 3
$ws.Range("D3").Value = '2.671.24'
$ws.Range("E3").Value = '  +0.13%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.40'
$ws.Range("E5").Value = '  -2.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.97'
$ws.Range("E6").Value = '  -3.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.522'
$ws.Range("E8").Value = '  -1.33%  '

# Row 9
$ws.Range("D9").Value = '2.670.02'
$ws.Range("E9").Value = '  +0.14%  '

# Row 10
$ws.Range("E10").Value = '  -5.96%  '

# Row 11
$ws.Range("E11").Value = '  +1.93%  '

# Row 12
$ws.Range("E12").Value = '  +0.32%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.98'
$ws.Range("E13").Value = '  -2.53%  '

# Row 14
$ws.Range("D14").Value = '3.170.52'
$ws.Range("E14").Value = '  +0.74%  '

# Row 15
$ws.Range("E15").Value = '  -5.14%  '

# Row 16
$ws.Range("D16").Value = '71.691.29'
$ws.Range("E16").Value = '  -1.43%  '

# Row 17
$ws.Range("E17").Value = '  -3.45%  '

# Row 18
$ws.Range("D18").Value = '2.669.88'
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  +4.68%  '

# Row 20
$ws.Range("E20").Value = '  +3.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.67'
$ws.Range("E21").Value = '  -3.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.17'
$ws.Range("E22").Value = '  -1.45%  '

# Row 23
$ws.Range("E23").Value = '  -1.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.89'
$ws.Range("E24").Value = '  -2.13%  '

# Row 25
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("E26").Value = '  -3.41%  '

# Row 27
$ws.Range("E27").Value = '  -2.77%  '

# Row 28
$ws.Range("D28").Value = '2.816.68'
$ws.Range("E28").Value = '  +0.36%  '

# Row 29
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0965'
$ws.Range("E30").Value = '  -1.73%  '

# Row 31
$ws.Range("E31").Value = '  -1.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '501.62'
$ws.Range("E32").Value = '  -8.45%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("E33").Value = '  -4.02%  '

# Row 34
$ws.Range("E34").Value = '  -2.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.42'
$ws.Range("E36").Value = '  -1.23%  '

# Row 37
$ws.Range("E37").Value = '  +0.26%  '

# Row 38
$ws.Range("E38").Value = '  -0.42%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.37'
$ws.Range("E39").Value = '  -3.72%  '

# Row 40
$ws.Range("E40").Value = '  -4.56%  '

# Row 41
$ws.Range("E41").Value = '  -5.26%  '

# Row 42
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("E43").Value = '  -3.02%  '

# Row 44
$ws.Range("E44").Value = '  -2.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.332'
$ws.Range("E45").Value = '  -1.58%  '

# Row 46
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.47'
$ws.Range("E46").Value = '  -0.66%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.66'
$ws.Range("E47").Value = '  +1.28%  '

# Row 48
$ws.Range("E48").Value = '  +2.11%  '

# Row 49
$ws.Range("E49").Value = '  -0.19%  '

# Row 50
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0754'
$ws.Range("E51").Value = '  -1.61%  '
